# References.xlsx edit: add a new reference row (Bennett & Brassard, BB84 /
# quantum cryptography paper) plus a new "Link" column with a hyperlink to
# the DOI, matching the commit "added python file, added a reference to
# References.xlsx".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fill in the new reference row's data -------------------------------
$ws.Range("B2").Value = "Charles H. Bennett, Gilles Brassard,"
$ws.Range("G1").Value = "Link"
$ws.Range("C2").Value = "Quantum cryptography: Public key distribution and coin tossing, Theoretical Computer Science,"

# Title cell gets vertically centered + wrapped text (long title).
$ws.Range("C2").VerticalAlignment = -4108
$ws.Range("C2").WrapText = $true

# Wrap text across the data rows so long descriptions/topics read cleanly
# (header row keeps the default style).
$ws.Range("A2:G16").WrapText = $true

# Hyperlink for the DOI link (auto-fills the cell text + "Hyperlink" style).
$ws.Hyperlinks.Add($ws.Range("G2"), "https://doi.org/10.1016/j.tcs.2014.05.025")

$ws.Range("E2").Value = "Quantum cryptography, QKD, BB84 Protocol"
$ws.Range("F2").Value = "Quantum cryptography, QKD, BB84 Protocol"
$ws.Range("D2").Value = 2014

# --- grow the table to cover the new column + extra blank rows ---------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G16"))
$tbl.ListColumns.Item(7).Name = "Link"

# --- column widths (widened to fit the new content) ---------------------
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth + 5
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(3).ColumnWidth + 10
$ws.Columns.Item(7).ColumnWidth = 34.24

# --- row height for the new, wrapped reference row -----------------------
$ws.Rows.Item(2).RowHeight = 34

# --- restore the active selection -----------------------------------
$ws.Range("F15").Select()
